$wb = $excel.ActiveWorkbook

$wsRes = $wb.Worksheets.Item("Resolutions")

# New data row 17 (mirrors the pattern used by existing rows, e.g. row 16)
$wsRes.Range("B17").Value = 1170
$wsRes.Range("C17").Value = 2532
$wsRes.Range("D17").Formula = "=B17/C17"
$wsRes.Range("E17").Formula = "=MIN(C17, B17)"
$wsRes.Range("F17").Value = 100
$wsRes.Range("G17").Formula = "=D17*F17"
$wsRes.Range("H17").Value = 0.1
$wsRes.Range("I17").Formula = "=(E17*H17) / 100"
$wsRes.Range("J17").Value = 1
$wsRes.Range("K17").Formula = "=F17*I17*J17"
$wsRes.Range("L17").Value = 0.1
$wsRes.Range("M17").Formula = "=B17*L17"
$wsRes.Range("N17").Formula = "=C17*L17"

# New data row 18
$wsRes.Range("B18").Value = 1848
$wsRes.Range("C18").Value = 2960
$wsRes.Range("D18").Formula = "=B18/C18"
$wsRes.Range("E18").Formula = "=MIN(C18, B18)"
$wsRes.Range("F18").Value = 100
$wsRes.Range("G18").Formula = "=D18*F18"
$wsRes.Range("H18").Value = 0.1
$wsRes.Range("I18").Formula = "=(E18*H18) / 100"
$wsRes.Range("J18").Value = 1
$wsRes.Range("K18").Formula = "=F18*I18*J18"
$wsRes.Range("L18").Value = 0.1
$wsRes.Range("M18").Formula = "=B18*L18"
$wsRes.Range("N18").Formula = "=C18*L18"

# Update the on-screen selections to match the authored view state.
$wsSample = $wb.Worksheets.Item("Sample-Resolutions")
$wsSample.Activate()
$wsSample.Range("C19:D19").Select()

$wsRes.Activate()
$wsRes.Range("M19").Select()
